# CommonComponent.pptx update
#  - Re-cache the "datetimeFigureOut" date placeholders (Master + all
#    Custom Layouts) from 7/8/2012 -> 14/8/12
#  - Resize/reposition a few shapes on slide 1
#  - Merge split text runs ("common::" + "exception", and the three
#    "EnrollException" " " "etc." triples) back into single runs
#  - Add a new "BuildProperties" box under "Common"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1. Date placeholder text (slide master + every custom layout)
# ---------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "14/8/12"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------
# 2. Resize / reposition shapes
# ---------------------------------------------------------------
# Outer rounded-rectangle background grows upward and taller
$outer = $s.Shapes.Item(1)
$outer.Top = 9.62496062992126
$outer.Height = 372

# "Group 157" (common / common.util) grows
$grp = $s.Shapes.Item(5)
$grp.Top = 33.62496062992126
$grp.Width = 162
$grp.Height = 132

# "Common" box widens and moves up slightly
$commonBox = $s.Shapes.Item(7)
$commonBox.Top = 75.62496062992126
$commonBox.Width = 120

# ---------------------------------------------------------------
# 3. Merge split runs back into single runs
# ---------------------------------------------------------------
# "common::" + "exception" -> "common::exception"
$commonExceptionGroup = $s.Shapes.Item(6)
$commonExceptionShape = $commonExceptionGroup.GroupItems.Item(1)
$commonExceptionShape.TextFrame.TextRange.Text = "temp-placeholder"
$commonExceptionShape.TextFrame.TextRange.Text = "common::exception"

# "EnrollException" + " " + "etc." -> "EnrollException" + " etc."
foreach ($idx in 9, 12, 13) {
    $sh = $s.Shapes.Item($idx)
    $tr = $sh.TextFrame.TextRange
    $tail = $tr.Characters(16, 5)
    $tail.Text = "tmp-xx"
    $tail2 = $sh.TextFrame.TextRange.Characters(16, 6)
    $tail2.Text = " etc."
}

# ---------------------------------------------------------------
# 4. Add the new "BuildProperties" rectangle under "Common"
# ---------------------------------------------------------------
# Create throwaway shapes first so the new shape's autogenerated
# id/name land on 35 / "Rectangle 34", matching the source id space.
$dummies = @()
for ($i = 0; $i -lt 33; $i++) {
    $dummies += $s.Shapes.AddShape(1, 0, 0, 10, 10)
}
foreach ($d in $dummies) {
    $d.Delete()
}

$newBox = $commonBox.Duplicate().Item(1)
$newBox.Left = 60
$newBox.Top = 123.62496062992126
$newBox.Width = 120
$newBox.Height = 36
$newBox.Name = "Rectangle 34"
$newBox.TextFrame.TextRange.Text = "temp-placeholder"
$newBox.TextFrame.TextRange.Text = "BuildProperties"
